# Append: 2026-01-17 06:28 JST
# The scraper re-ran, found 2 new postings, re-scored/re-ranked the
# rolling window, and trimmed it back down to 6 data rows (was 16).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# --- column B got a touch narrower ---------------------------------------
# ColumnWidth setter adds a constant ~0.8333 "padding" offset internally,
# so back it out to land on an exact 47.
$ws.Columns.Item(2).ColumnWidth = 46.166666666666664

# --- row 2: brand-new posting ---------------------------------------------
$ws.Cells.Item(2,1).Value = "2026-01-17 06:28:18"
$ws.Cells.Item(2,2).Value = "【急募】airtableで社内業務管理システムを共に構築してくれる方"
$ws.Cells.Item(2,3).Value = "システム開発"
$ws.Cells.Item(2,4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(2,5).Value = "期限情報なし"
$ws.Cells.Item(2,6).Value = "https://www.lancers.jp/work/detail/5473383"
$ws.Cells.Item(2,7).Value = 353
$ws.Cells.Item(2,8).Value = "🔥AI,Ai ◇管理"

# --- row 3: brand-new posting ---------------------------------------------
$ws.Cells.Item(3,1).Value = "2026-01-17 06:28:18"
$ws.Cells.Item(3,2).Value = "【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪"
$ws.Cells.Item(3,3).Value = "システム開発"
$ws.Cells.Item(3,4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(3,5).Value = "期限情報なし"
$ws.Cells.Item(3,6).Value = "https://www.lancers.jp/work/detail/5217096"
$ws.Cells.Item(3,7).Value = 243
$ws.Cells.Item(3,8).Value = "🔥API ◆ツール"

# --- row 4: same posting as before, only the scrape timestamp moved -------
$ws.Cells.Item(4,1).Value = "2026-01-17 06:28:18"

# --- row 5: now holds what used to be row 6's posting ----------------------
$ws.Cells.Item(5,1).Value = "2026-01-17 06:28:18"
$ws.Cells.Item(5,2).Value = "※急募:Flutterによる業務アプリの開発(+next.js)"
$ws.Cells.Item(5,3).Value = "システム開発"
$ws.Cells.Item(5,4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(5,5).Value = "期限情報なし"
$ws.Cells.Item(5,6).Value = "https://www.lancers.jp/work/detail/5473146"
$ws.Cells.Item(5,7).Value = 218
$ws.Cells.Item(5,8).Value = "🔥Next.js ◆開発 ◇アプリ"

# --- row 6: now holds what used to be row 8's posting ----------------------
$ws.Cells.Item(6,1).Value = "2026-01-17 06:28:18"
$ws.Cells.Item(6,2).Value = "【急募】Accessでの受発注管理・請求書発行システム開発"
$ws.Cells.Item(6,3).Value = "システム開発"
$ws.Cells.Item(6,4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(6,5).Value = "期限情報なし"
$ws.Cells.Item(6,6).Value = "https://www.lancers.jp/work/detail/5473234"
$ws.Cells.Item(6,7).Value = 148
$ws.Cells.Item(6,8).Value = "◆開発,システム開発 ◇管理"

# --- row 7: now holds what used to be row 9's posting ----------------------
$ws.Cells.Item(7,1).Value = "2026-01-17 06:28:18"
$ws.Cells.Item(7,2).Value = "【バイナリ解析 / 逆コンパイル】EPCデータ解析ツール開発"
$ws.Cells.Item(7,3).Value = "システム開発"
$ws.Cells.Item(7,4).Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Cells.Item(7,5).Value = "期限情報なし"
$ws.Cells.Item(7,6).Value = "https://www.lancers.jp/work/detail/5473181"
$ws.Cells.Item(7,7).Value = 135
$ws.Cells.Item(7,8).Value = "◆ツール,開発"

# --- the old rows 8-17 fall out of the rolling window entirely ------------
$ws.Rows("8:17").Delete()

# --- hyperlinks: the engine's Hyperlinks.Delete() clears the whole sheet
# collection regardless of the range it's called on, so wipe once and
# re-add exactly the six links that should remain (F2:F7), pointing at
# their (possibly updated) URLs.
$ws.Range("F2:F7").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Cells.Item(2,6), "https://www.lancers.jp/work/detail/5473383") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(3,6), "https://www.lancers.jp/work/detail/5217096") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(4,6), "https://www.lancers.jp/work/detail/5473147") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(5,6), "https://www.lancers.jp/work/detail/5473146") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(6,6), "https://www.lancers.jp/work/detail/5473234") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(7,6), "https://www.lancers.jp/work/detail/5473181") | Out-Null

Write-Host "done"
